$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.796.07'
$ws.Range("E2").Value = '  -0.76%  '

$ws.Range("D3").Value = '2.322.83'
$ws.Range("E3").Value = '  -0.26%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '301.85'
$ws.Range("E5").Value = '  -0.95%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.94'
$ws.Range("E6").Value = '  -3.96%  '

$ws.Range("E7").Value = '  -0.73%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.494'
$ws.Range("E9").Value = '  -1.95%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.96'
$ws.Range("E10").Value = '  -4.64%  '

$ws.Range("E11").Value = '  -2.20%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.73'
$ws.Range("E12").Value = '  -3.59%  '

$ws.Range("E13").Value = '  +1.92%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.71'
$ws.Range("E14").Value = '  -3.18%  '

$ws.Range("D15").Value = '2.684.99'
$ws.Range("E15").Value = '  -0.28%  '

$ws.Range("D16").Value = '2.338.32'
$ws.Range("E16").Value = '  +0.63%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.788'
$ws.Range("E17").Value = '  -0.02%  '

$ws.Range("D18").Value = '42.740.59'
$ws.Range("E18").Value = '  -0.67%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.02'
$ws.Range("E19").Value = '  -4.11%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.17'
$ws.Range("E20").Value = '  +1.45%  '

$ws.Range("E21").Value = '  -1.69%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.86'
$ws.Range("E22").Value = '  -0.21%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.39'
$ws.Range("E23").Value = '  -0.91%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.22'
$ws.Range("E24").Value = '  +0.93%  '

$ws.Range("E25").Value = '  +0.03%  '

$ws.Range("E26").Value = '  -1.39%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.54'
$ws.Range("E27").Value = '  -1.53%  '

$ws.Range("E28").Value = '  -1.09%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.09'
$ws.Range("E29").Value = '  -0.49%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.41'
$ws.Range("E30").Value = '  -5.37%  '

$ws.Range("E31").Value = '  +0.03%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '139.88'
$ws.Range("E32").Value = '  -15.95%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.99'
$ws.Range("E33").Value = '  -0.44%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.45'
$ws.Range("E34").Value = '  -3.07%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0697'
$ws.Range("E35").Value = '  +0.06%  '

$ws.Range("E36").Value = '  -0.86%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.35'
$ws.Range("E37").Value = '  -4.61%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.81'
$ws.Range("E38").Value = '  +2.61%  '

$ws.Range("E39").Value = '  -0.56%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.52'
$ws.Range("E40").Value = '  +24.37%  '

$ws.Range("E41").Value = '  -2.36%  '

$ws.Range("E42").Value = '  -1.10%  '

$ws.Range("D43").Value = '1.934.87'
$ws.Range("E43").Value = '  -3.16%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0279'
$ws.Range("E44").Value = '  -0.74%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.24'
$ws.Range("E45").Value = '  -4.54%  '

$ws.Range("E46").Value = '  -0.74%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.71'
$ws.Range("E47").Value = '  -2.70%  '

$ws.Range("E48").Value = '  +0.32%  '

$ws.Range("D49").Value = '2.552.02'
$ws.Range("E49").Value = '  -0.25%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '52.54'
$ws.Range("E50").Value = '  -2.28%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.09'
$ws.Range("E51").Value = '  +0.16%  '
